$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '42.674.03'
$ws.Range('E2').Value = '  +2.10%  '
Set-TextValue 'D3' '2.286.38'
$ws.Range('E3').Value = '  +3.58%  '
$ws.Range('E4').Value = '  +0.12%  '
Set-TextValue 'D5' '251.04'
$ws.Range('E5').Value = '  +0.09%  '
Set-TextValue 'D6' '0.634'
$ws.Range('E6').Value = '  +2.58%  '
Set-TextValue 'D7' '72.61'
$ws.Range('E7').Value = '  +6.49%  '
$ws.Range('E8').Value = '  +0.03%  '
Set-TextValue 'D9' '0.647'
$ws.Range('E9').Value = '  +2.50%  '
Set-TextValue 'D10' '39.28'
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D11' '59.11'
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D12' '0.0961'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('E13').Value = '  +3.75%  '
$ws.Range('E14').Value = '  +2.05%  '
Set-TextValue 'D15' '2.630.68'
$ws.Range('E15').Value = '  +3.72%  '
Set-TextValue 'D16' '15.15'
$ws.Range('E16').Value = '  +3.74%  '
Set-TextValue 'D17' '0.882'
$ws.Range('E17').Value = '  +1.49%  '
Set-TextValue 'D18' '2.280.42'
$ws.Range('E18').Value = '  +4.06%  '
Set-TextValue 'D19' '42.641.19'
$ws.Range('E19').Value = '  +2.27%  '
Set-TextValue 'D20' '0.0₃0998'
$ws.Range('E20').Value = '  +3.98%  '
$ws.Range('E21').Value = '  +1.76%  '
Set-TextValue 'D22' '72.42'
$ws.Range('E22').Value = '  +0.18%  '
Set-TextValue 'D23' '2.26'
$ws.Range('E23').Value = '  +9.75%  '
Set-TextValue 'D24' '235.50'
$ws.Range('E24').Value = '  +1.63%  '
Set-TextValue 'D25' '3.94'
$ws.Range('E25').Value = '  +0.68%  '
Set-TextValue 'D26' '11.60'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('E30').Value = '  +2.42%  '
Set-TextValue 'D31' '167.22'
$ws.Range('E31').Value = '  -0.08%  '
Set-TextValue 'D32' '21.13'
$ws.Range('E32').Value = '  +3.74%  '
Set-TextValue 'D33' '6.46'
$ws.Range('E33').Value = '  +10.27%  '
Set-TextValue 'D34' '0.127'
$ws.Range('E34').Value = '  +4.87%  '
Set-TextValue 'D35' '0.0810'
$ws.Range('E35').Value = '  +2.20%  '
Set-TextValue 'D36' '30.95'
$ws.Range('E36').Value = '  +18.20%  '
$ws.Range('E37').Value = '  +2.77%  '
Set-TextValue 'D38' '4.75'
$ws.Range('E38').Value = '  +13.98%  '
Set-TextValue 'D39' '4.75'
$ws.Range('E39').Value = '  +3.08%  '
Set-TextValue 'D40' '0.0308'
$ws.Range('E40').Value = '  -0.62%  '
Set-TextValue 'D41' '14.14'
$ws.Range('E41').Value = '  +15.89%  '
$ws.Range('E42').Value = '  +3.95%  '
Set-TextValue 'D43' '5.95'
$ws.Range('E43').Value = '  +5.61%  '
Set-TextValue 'D44' '0.216'
$ws.Range('E44').Value = '  +9.25%  '
Set-TextValue 'D45' '9.21'
$ws.Range('E45').Value = '  +7.53%  '
Set-TextValue 'D46' '61.94'
$ws.Range('E46').Value = '  -0.46%  '
Set-TextValue 'D47' '4.89'
$ws.Range('E47').Value = '  -5.23%  '
Set-TextValue 'D48' '0.104'
$ws.Range('E48').Value = '  +3.66%  '
Set-TextValue 'D49' '1.19'
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('E50').Value = '  +0.19%  '
Set-TextValue 'D51' '97.24'
$ws.Range('E51').Value = '  +5.86%  '
